# preparation publication 0.2.0
# - Add a new "Jurisdiction" / "iso:code:3166:FR" property row to the
#   Metadata sheet (inserted right after the "Contact" row and before
#   "Description", pushing all subsequent property rows down by one).
# - Bump the Version property from 0.1.1 to 0.2.0.
# - Refresh the Date property to the new publication timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new blank row before the current row 11 ("Description"),
# shifting everything from row 11 downward to row 12 onward.
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row with the Jurisdiction property.
$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = "iso:code:3166:FR"

# Match the formatting used by the rest of the property table (the row
# below, which used to be row 11 and is now row 12, has the correct
# body style) instead of the blank default style Insert() applied.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the Version property value.
$ws.Cells.Item(3, 2).Value = "0.2.0"

# Update the Date property value.
$ws.Cells.Item(8, 2).Value = "2023-10-20T08:59:58+00:00"
